# Generate Report for Handback
# Update the "generated" timestamps on the Overview / zh-cn / de-de sheets
# to reflect a fresh report run.

$wb = $excel.ActiveWorkbook

# --- Overview sheet -------------------------------------------------
# "Latest HO Xliff Generate Date" (G2): 2016-09-05 17:14:28 -> 17:15:43
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-09-05 17:15:43"

# --- zh-cn sheet ------------------------------------------------------
# "Correspond Handoff Datetime" (H2): 2016-09-05 17:14:23 -> 17:15:38
# "Correspond Handback DateTime" (K2): 2016-09-05 17:15:11 -> 17:15:56
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-09-05 17:15:38"
$wsZhCn.Range("K2").Value = "2016-09-05 17:15:56"

# --- de-de sheet ------------------------------------------------------
# "Correspond Handback DateTime" (K2): 2016-09-05 17:15:22 -> 17:16:10
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("K2").Value = "2016-09-05 17:16:10"
